# Applies the edits described in the diff:
# - updates the date heading
# - updates each division-problem cell text in the worksheet table
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-11 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-12 Tuesday", 2)
$d.Content.Find.Execute("55÷3=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "27÷2=13, 1", 2)
$d.Content.Find.Execute("57÷5=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "86÷6=14, 2", 2)
$d.Content.Find.Execute("28÷5=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "48÷6=8, 0", 2)
$d.Content.Find.Execute("23÷8=2, 7", $true, $false, $false, $false, $false, $true, 1, $false, "64÷9=7, 1", 2)
$d.Content.Find.Execute("13÷5=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "87÷6=14, 3", 2)
$d.Content.Find.Execute("71÷7=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "60÷7=8, 4", 2)
$d.Content.Find.Execute("65÷2=32, 1", $true, $false, $false, $false, $false, $true, 1, $false, "86÷6=14, 2", 2)
$d.Content.Find.Execute("38÷3=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "64÷3=21, 1", 2)
$d.Content.Find.Execute("28÷4=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "69÷7=9, 6", 2)
$d.Content.Find.Execute("80÷2=40, 0", $true, $false, $false, $false, $false, $true, 1, $false, "20÷8=2, 4", 2)
$d.Content.Find.Execute("17÷4=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "41÷2=20, 1", 2)
$d.Content.Find.Execute("84÷2=42, 0", $true, $false, $false, $false, $false, $true, 1, $false, "68÷5=13, 3", 2)
$d.Content.Find.Execute("19÷2=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "43÷7=6, 1", 2)
$d.Content.Find.Execute("13÷6=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "46÷7=6, 4", 2)
$d.Content.Find.Execute("27÷6=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "71÷4=17, 3", 2)
$d.Content.Find.Execute("98÷7=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "78÷5=15, 3", 2)
$d.Content.Find.Execute("20÷5=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "96÷4=24, 0", 2)
$d.Content.Find.Execute("42÷9=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "83÷2=41, 1", 2)
$d.Content.Find.Execute("19÷6=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "37÷8=4, 5", 2)
$d.Content.Find.Execute("53÷8=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "26÷8=3, 2", 2)
$d.Content.Find.Execute("89÷2=44, 1", $true, $false, $false, $false, $false, $true, 1, $false, "53÷9=5, 8", 2)
$d.Content.Find.Execute("13÷7=1, 6", $true, $false, $false, $false, $false, $true, 1, $false, "73÷3=24, 1", 2)
$d.Content.Find.Execute("94÷8=11, 6", $true, $false, $false, $false, $false, $true, 1, $false, "33÷7=4, 5", 2)
$d.Content.Find.Execute("14÷2=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "76÷7=10, 6", 2)
$d.Content.Find.Execute("62÷8=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "79÷9=8, 7", 2)
